# Scrape time randomized to mimic human behavior.
# Rename the three sheets to tag them with the scraped ticker (RDFN),
# and bump the recorded scrape timestamps forward (new scrape run).

$wb = $excel.ActiveWorkbook

$wsTechnicals = $wb.Worksheets.Item("technicals")
$wsPuts       = $wb.Worksheets.Item("puts")
$wsCalls      = $wb.Worksheets.Item("calls")

$wsTechnicals.Name = "technicalsRDFN"
$wsPuts.Name       = "putsRDFN"
$wsCalls.Name      = "callsRDFN"

# technicals: single curr_time cell at K2
$wsTechnicals.Range("K2").Value = "2020-02-20 22:04:50"

# puts: curr_time column D, rows 2-24
for ($row = 2; $row -le 24; $row++) {
    $wsPuts.Cells.Item($row, 4).Value = "2020-02-20 22:04:51"
}

# calls: curr_time column D, rows 2-27
for ($row = 2; $row -le 27; $row++) {
    $wsCalls.Cells.Item($row, 4).Value = "2020-02-20 22:04:53"
}
